# Apply the "Latest updates from master" revision to sample_inputs.xlsx
# (AAE 560 Air Taxi Model) via Excel COM automation.

$wb = $excel.ActiveWorkbook

$wsOperator = $wb.Worksheets.Item("Operator")
$wsAircraft = $wb.Worksheets.Item("Aircraft")
$wsPorts    = $wb.Worksheets.Item("Ports")

# ---------------------------------------------------------------------
# Aircraft sheet: bump the Speed (mph) figures for both aircraft types
# ---------------------------------------------------------------------
$wsAircraft.Range("B2").Value = 170
$wsAircraft.Range("B3").Value = 170

# ---------------------------------------------------------------------
# Ports sheet: Landing Slots column (E) goes from 1 to 5 for every
# existing port, and a brand-new Port 6 row is appended.
# ---------------------------------------------------------------------
$wsPorts.Range("E2").Value = 5
$wsPorts.Range("E3").Value = 5
$wsPorts.Range("E4").Value = 5
$wsPorts.Range("E5").Value = 5
$wsPorts.Range("E6").Value = 5

# Bring row 7 formatting in line with row 6 before filling it in so the
# percentage/number styles (cols D & F) carry over to the new row.
$wsPorts.Range("A6:F6").Copy()
$wsPorts.Range("A7:F7").PasteSpecial(-4122)

$wsPorts.Range("A7").Value = 6
$wsPorts.Range("B7").Value = 50
$wsPorts.Range("C7").Value = 72
$wsPorts.Range("D7").Value = 0.6
$wsPorts.Range("E7").Value = 5
$wsPorts.Range("F7").Value = 20

# The now-unused helper cell I1 (stray formatting only, no data) is
# cleared so the sheet's used range shrinks back down to F.
$wsPorts.Range("I1").Clear()

# ---------------------------------------------------------------------
# Operator sheet: fleet mix, charging/port selections, and the
# helper table that now spans an extra "Port 6" column (H).
# ---------------------------------------------------------------------

# Fleet size: 2/1 Type1/Type2 aircraft -> 8/8
$wsOperator.Range("C10").Value = 8
$wsOperator.Range("D10").Value = 8

# Serviced-ports row: G13 flips from "No" to "Yes", and the newly
# visible H13 cell (Port 6) is also serviced ("Yes"), matching the
# input-cell styling used by the rest of the row.
$wsOperator.Range("G13").Value = "Yes"
$wsOperator.Range("G13").Copy()
$wsOperator.Range("H13").PasteSpecial(-4122)
$wsOperator.Range("H13").Value = "Yes"

# Charging-equipment row: H14 (Port 6) gets a "Slow" charger, again
# copying the formatting used by the sibling cells in the row.
$wsOperator.Range("G14").Copy()
$wsOperator.Range("H14").PasteSpecial(-4122)
$wsOperator.Range("H14").Value = "Slow"

# Serviced-ports count formula now covers the extended C:I range.
$wsOperator.Range("B13").Formula = '=COUNTIF(C13:I13,"Yes")'

# Charger-type data validation list extends from C14:G14 to C14:H14.
$wsOperator.Range("C14:G14").Validation.Delete()
$wsOperator.Range("C14:H14").Validation.Add(3, 1, 1, '"Slow, Fast, None"')

# ---------------------------------------------------------------------
# View state: Ports is now the active/selected sheet.
# ---------------------------------------------------------------------
$wsOperator.Range("A1").Select() | Out-Null
$wsAircraft.Range("A1").Select() | Out-Null
$wsPorts.Range("A1").Select() | Out-Null
$wsPorts.Activate() | Out-Null

$wb.Application.Calculate() | Out-Null
